$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column A (xlsx stored width 16.42578125 -> 15.42578125).
# The runtime quantizes ColumnWidth to 1/6-character steps, so we feed it
# the value whose pt*2 px rounding lands on the bucket closest to the target.
$ws.Columns.Item(1).ColumnWidth = 14.592447916666666

# Replace the 33 sampled values in A1:A33 with the new dataset
$values = @(0.21888754210301897, -0.0059999999754793976, -0.0039999999781912265, -0.0079999999599902338, -0.002999999977112644, -0.0019999999747444264, -0.0099999999434996489, -0.0099999999447746291, -0.0019999999789255263, -0.0019999999818249847, -0.0029999999782459597, 0.046241383756871368, -0.0034999999774871782, -0.007999999959966253, 0.021146649390370698, -0.0019999999835818016, -0.0019999999831101789, -0.0039999999751332282, -0.0039999999828941313, -0.039781250589815187, -0.0039999999778714823, -0.0039999999777320383, -0.0049999999729459788, -0.019999999911735955, -0.019999999910611521, -0.0024999999713273269, -0.020250536675115782, -0.0019999999631163945, -0.0069999999376495481, -0.059999999726513931, -0.0069999999329155571, -0.0099999999208275625, -0.0039999999443285361)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
